# edit.ps1 - apply the changes described by the diff
#
# 1) "Leggere stato p" + hidden _GoBack bookmark + "orta (aperta o chiusa)."
#    -> merge into a single run "Leggere stato porta (aperta o chiusa)."
#    (this also removes the old _GoBack bookmark, since it sat inside the
#    replaced range)
# 2) Insert a new, empty paragraph (with its own formatting, carrying the
#    _GoBack bookmark) right before the "Strumenti e metodi" paragraph.
# 3) Move <w:lastRenderedPageBreak/> from the "Elenco dei principali..."
#    paragraph to the "Strumenti e metodi" paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: fix the split run / stray bookmark around "Leggere stato porta"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Leggere stato porta (aperta o chiusa).", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Leggere stato porta (aperta o chiusa).", 2) | Out-Null

# ---------------------------------------------------------------------
# Step 2: insert the new empty paragraph before "Strumenti e metodi"
# ---------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Porta chiusa.`r") {
        $targetIdx = $i
        break
    }
}

$portaChiusa = $d.Paragraphs.Item($targetIdx)
$portaChiusa.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item($targetIdx + 1)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1854"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="160"/><w:ind w:left="142"/><w:contextualSpacing/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:lang w:val="it-CH" w:eastAsia="it-CH"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# Step 3: move <w:lastRenderedPageBreak/> from the "Elenco dei
# principali..." paragraph onto the "Strumenti e metodi" paragraph
# ---------------------------------------------------------------------
$strumentiIdx = 0
$elencoIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Strumenti e metodi`r") {
        $strumentiIdx = $i
    } elseif ($t -like "Elenco dei principali strumenti e metodi*") {
        $elencoIdx = $i
    }
}

$strumenti = $d.Paragraphs.Item($strumentiIdx)
$strumentiXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Strumenti e metodi</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$strumenti.Range.InsertXML($strumentiXml)

$elenco = $d.Paragraphs.Item($elencoIdx)
$elencoXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:i/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:i/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Elenco dei principali strumenti e metodi (hardware, software, linguaggio di programmazione, etc.) che devono essere utilizzati per la soluzione</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$elenco.Range.InsertXML($elencoXml)

Write-Output "edit complete"
